$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "checker graphics" task duration reduced to 0 (implemented counter/checker graphics)
$ws.Range("C4").Value = 0

# "check that" task duration reduced to 0 (added bug desc.)
$ws.Range("C5").Value = 0

# Move the active selection from B10 to F10
$ws.Range("F10").Select() | Out-Null
